$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 15, shifting existing rows 15-65 down to 16-66.
$ws.Rows("15:15").Insert()

# Populate the new row 15 with the new record's data.
$ws.Range("A15").Value = 5
$ws.Range("B15").Value = "Macroferia Regional de Talca"
$ws.Range("C15").Value = "Maule"
$ws.Range("D15").Value2 = 44487
$ws.Range("E15").Value = 7
$ws.Range("F15").Value = 100112013
$ws.Range("G15").Value = "Alcachofa"
$ws.Range("H15").Value = "Madrigal"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 400
$ws.Range("K15").Value = 9000
$ws.Range("L15").Value = 9000
$ws.Range("M15").Value = 9000
$ws.Range("N15").Value = "$/caja 40 unidades"
$ws.Range("O15").Value = "Región del Maule"
$ws.Range("P15").Value = 225
$ws.Range("Q15").Value = 40
$ws.Range("R15").Value = "Hortaliza"
